# Update generator linear examples: refresh the randomly-generated
# numeric samples on the "follower restrictions", "modified point",
# "vector bf" and "vector BF" sheets.
#
# NOTE: cells in this workbook store their numbers as text (shared
# strings, t="s"), not as real numbers. Assigning a numeric-looking
# string straight to .Value would make Excel auto-coerce it into a
# number cell. To keep the cells as text we temporarily force a text
# number format ("@") before writing the value, then restore the
# original ("Normal") style afterwards.

$wb = $excel.ActiveWorkbook

# --- Restricciones_del_follower -------------------------------------
$ws = $wb.Worksheets.Item("Restricciones_del_follower")
$rng = $ws.Range("B2:F5")
$rng.NumberFormat = "@"

$ws.Range("B2").Value = "-4.374623078112156"
$ws.Range("D2").Value = "0.33468162538227564"
$ws.Range("E2").Value = "0"
$ws.Range("F2").Value = "0.9872811638058144"

$ws.Range("B3").Value = "-2.8102693382873367"
$ws.Range("D3").Value = "0.9092567913461869"
$ws.Range("F3").Value = "0.2538126765657339"

$ws.Range("B4").Value = "0.9341385726238034"
$ws.Range("D4").Value = "0.7906785535517057"
$ws.Range("E4").Value = "0.8266745873534492"

$ws.Range("B5").Value = "0.36494658748581443"
$ws.Range("D5").Value = "0.5618257705012442"
$ws.Range("E5").Value = "0.3714762402532654"
$ws.Range("F5").Value = "0.4485872944177247"

$rng.Style = "Normal"

# --- Punto_modificado -------------------------------------------------
$ws = $wb.Worksheets.Item("Punto_modificado")
$rng = $ws.Range("A2:B2")
$rng.NumberFormat = "@"

$ws.Range("A2").Value = "4.184892416399492"
$ws.Range("B2").Value = "4.374623078112156"

$rng.Style = "Normal"

# --- Vector_bf -----------------------------------------------------------
# NOTE: worksheet name lookup is case-insensitive, and this workbook has
# both "Vector_bf" (sheet 5) and "Vector_BF" (sheet 6), so they are
# addressed by their 1-based sheet index instead of by name to avoid
# ambiguity.
$ws = $wb.Worksheets.Item(5)
$rng = $ws.Range("A2:A2")
$rng.NumberFormat = "@"

$ws.Range("A2").Value = "-2.5941065025660786"

$rng.Style = "Normal"

# --- Vector_BF ------------------------------------------------------------
$ws = $wb.Worksheets.Item(6)
$rng = $ws.Range("A2:A3")
$rng.NumberFormat = "@"

$ws.Range("A2").Value = "-1.3125795483665108"
$ws.Range("A3").Value = "1.718127065546367"

$rng.Style = "Normal"

Write-Host "Generator values updated"
